$d = $word.ActiveDocument

# --- The "Table" style sample table -----------------------------------
$tbl = $d.Tables(1)

# Table now stretches to 100% of the available width (tblW -> 5000 pct)
# instead of a fixed 1532 dxa width.
$tbl.PreferredWidthType = 2   # wdPreferredWidthPercent
$tbl.PreferredWidth = 250     # engine scales *20 -> w:w="5000" w:type="pct"

# Column widths: 2093 dxa / 7483 dxa (applies to both the <w:tblGrid>
# and each row's <w:tcW>).
$tbl.Columns(1).Width = 104.65   # 2093 dxa / 20
$tbl.Columns(2).Width = 374.15   # 7483 dxa / 20

# Header row text: "Table" / "Table" -> "Column1" / "Column2LveryLongName..."
# (Find.Execute's "replace" ignores Range scoping in this runtime and hits
# the whole story, so set the cell Range.Text directly instead - it stays
# correctly confined to the cell.)
$tbl.Cell(1, 1).Range.Text = "Column1"
$tbl.Cell(1, 2).Range.Text = "Column2LveryLongName1234512345123451234512345"

# --- The "Table" table style -------------------------------------------
$style = $d.Styles("Table")
$style.BaseStyle = $d.Styles("TableNormal")
